$d = $word.ActiveDocument
$r = $d.Content
$anchor = "Ako klijent potražuje iznos koji je veći od njegove prosečne plate za poslednja 3 meseca, u zavisnosti od tipa računa maksimalno može biti iznos do 3 mesečne plate"
$ok = $r.Find.Execute($anchor, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pStart = $r.Start
$midStart = $pStart + 17
$midEnd = $pStart + 148
$midRange = $d.Range($midStart, $midEnd)
$searchText = $midRange.Text
Write-Output "searchText len=$($searchText.Length)"
$newMid = "žuje iznos veći od jedne mesečne plate za Premium i 3 "
$ok2 = $midRange.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, $newMid, 2)
Write-Output "ok2=$ok2 newStart=$($midRange.Start) newEnd=$($midRange.End)"
